$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H12").Value = 498
$ws.Range("I12").Value = 498
$ws.Range("K12").Value = 498
$ws.Range("M12").Value = -328
$ws.Range("H93").Value = 40601
$ws.Range("J93").Value = 40601
$ws.Range("L93").Value = 40601
$ws.Range("N93").Value = -45593
$ws.Range("H98").Value = 1863.421
$ws.Range("I98").Value = 1503.8
$ws.Range("J98").Value = 3212
$ws.Range("K98").Value = 1503.8
$ws.Range("L98").Value = 3212
$ws.Range("M98").Value = -5.799999999999955
$ws.Range("N98").Value = -6208
$ws.Range("H101").Value = 476.4
$ws.Range("I101").Value = 449.375
$ws.Range("K101").Value = 1348.125
$ws.Range("M101").Value = 273.875
$ws.Range("H106").Value = 8500
$ws.Range("I106").Value = 8500
$ws.Range("K106").Value = 8500
$ws.Range("M106").Value = -7869
$ws.Range("H113").Value = 3334.375
$ws.Range("I113").Value = 2981.6
$ws.Range("K113").Value = 2981.6
$ws.Range("M113").Value = 272.4000000000001
$ws.Range("H122").Value = 1863.421
$ws.Range("I122").Value = 1503.8
$ws.Range("J122").Value = 3212
$ws.Range("K122").Value = 4511.4
$ws.Range("L122").Value = 9636
$ws.Range("M122").Value = -2061.4
$ws.Range("N122").Value = -14536
$ws.Range("H137").Value = 2365.5833
$ws.Range("I137").Value = 1888.5555
$ws.Range("J137").Value = 3796.6667
$ws.Range("K137").Value = 5665.666499999999
$ws.Range("L137").Value = 11390.0001
$ws.Range("M137").Value = -3115.666499999999
$ws.Range("N137").Value = -16490.0001

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 3996.25
$ws.Range("I2").Value = 3996.25
$ws.Range("K2").Value = 3996.25
$ws.Range("M2").Value = -3883.25
$ws.Range("H63").Value = 1779.6
$ws.Range("I63").Value = 1166
$ws.Range("J63").Value = 2700
$ws.Range("K63").Value = 1166
$ws.Range("L63").Value = 2700
$ws.Range("M63").Value = -480
$ws.Range("N63").Value = -4072
$ws.Range("H66").Value = 1779.6
$ws.Range("I66").Value = 1166
$ws.Range("J66").Value = 2700
$ws.Range("K66").Value = 5830
$ws.Range("L66").Value = 13500
$ws.Range("M66").Value = -2398
$ws.Range("N66").Value = -20364
$ws.Range("H97").Value = 1301.6666
$ws.Range("H116").Value = 3996.25
$ws.Range("I116").Value = 3996.25
$ws.Range("K116").Value = 3996.25
$ws.Range("M116").Value = -1702.25
$ws.Range("H122").Value = 2449.75
$ws.Range("I122").Value = 2449.75
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7349.25
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -4899.25
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 920
$ws.Range("I132").Value = 920
$ws.Range("K132").Value = 2760
$ws.Range("M132").Value = -230

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 3996.25
$ws.Range("I3").Value = 3996.25
$ws.Range("K3").Value = 3996.25
$ws.Range("M3").Value = -3882.25
$ws.Range("H35").Value = 12666.667
$ws.Range("J35").Value = 12666.667
$ws.Range("L35").Value = 12666.667
$ws.Range("N35").Value = -13286.667
$ws.Range("H59").Value = 77890
$ws.Range("J59").Value = 77890
$ws.Range("L59").Value = 77890
$ws.Range("N59").Value = -79584
$ws.Range("H94").Value = 2932.6667
$ws.Range("I94").Value = 2932.6667
$ws.Range("K94").Value = 2932.6667
$ws.Range("M94").Value = -2481.6667
$ws.Range("H107").Value = 1322
$ws.Range("I107").Value = 1322
$ws.Range("K107").Value = 1322
$ws.Range("M107").Value = 598
$ws.Range("H110").Value = 59666.668
$ws.Range("J110").Value = 59666.668
$ws.Range("L110").Value = 59666.668
$ws.Range("N110").Value = -67846.66800000001
$ws.Range("H134").Value = 6054.9653
$ws.Range("I134").Value = 6054.9653
$ws.Range("K134").Value = 18164.8959
$ws.Range("M134").Value = -15629.8959

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H6").Value = 621.8
$ws.Range("I6").Value = 1277.75
$ws.Range("J6").Value = 184.5
$ws.Range("K6").Value = 1277.75
$ws.Range("L6").Value = 184.5
$ws.Range("M6").Value = -1164.75
$ws.Range("N6").Value = -410.5
$ws.Range("H16").Value = 6405.25
$ws.Range("I16").Value = 1873.6666
$ws.Range("K16").Value = 1873.6666
$ws.Range("M16").Value = -1586.6666
$ws.Range("H22").Value = 500
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 500
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1200
$ws.Range("H86").Value = 4601.2
$ws.Range("I86").Value = 3668.6667
$ws.Range("K86").Value = 3668.6667
$ws.Range("M86").Value = -2545.6667
$ws.Range("H89").Value = 4601.2
$ws.Range("I89").Value = 3668.6667
$ws.Range("K89").Value = 18343.3335
$ws.Range("M89").Value = -12727.3335
$ws.Range("H99").Value = 4750
$ws.Range("I99").Value = 4572.857
$ws.Range("J99").Value = 4998
$ws.Range("K99").Value = 4572.857
$ws.Range("L99").Value = 4998
$ws.Range("M99").Value = -3074.857
$ws.Range("N99").Value = -7994
$ws.Range("H113").Value = 6405.25
$ws.Range("I113").Value = 1873.6666
$ws.Range("K113").Value = 1873.6666
$ws.Range("M113").Value = 296.3334
$ws.Range("H122").Value = 1566.6666
$ws.Range("I122").Value = 1380
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 4140
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -1690
$ws.Range("N122").Value = -12400
$ws.Range("H126").Value = 4750
$ws.Range("I126").Value = 4572.857
$ws.Range("J126").Value = 4998
$ws.Range("K126").Value = 13718.571
$ws.Range("L126").Value = 14994
$ws.Range("M126").Value = -11248.571
$ws.Range("N126").Value = -19934
$ws.Range("H134").Value = 2433.2
$ws.Range("I134").Value = 2686.647
$ws.Range("K134").Value = 8059.941
$ws.Range("M134").Value = -5524.941

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H5").Value = 1447.4166
$ws.Range("I5").Value = 1639.4
$ws.Range("J5").Value = 487.5
$ws.Range("K5").Value = 4918.200000000001
$ws.Range("L5").Value = 1462.5
$ws.Range("M5").Value = -4806.200000000001
$ws.Range("N5").Value = -1686.5
$ws.Range("H81").Value = 5169.1665
$ws.Range("I81").Value = 5000
$ws.Range("J81").Value = 5203
$ws.Range("K81").Value = 15000
$ws.Range("L81").Value = 15609
$ws.Range("M81").Value = -13877
$ws.Range("N81").Value = -17855
$ws.Range("H84").Value = 5169.1665
$ws.Range("I84").Value = 5000
$ws.Range("J84").Value = 5203
$ws.Range("K84").Value = 45000
$ws.Range("L84").Value = 46827
$ws.Range("M84").Value = -39384
$ws.Range("N84").Value = -58059
$ws.Range("H113").Value = 1951.5
$ws.Range("I113").Value = 1503
$ws.Range("K113").Value = 4509
$ws.Range("M113").Value = -2339
$ws.Range("H135").Value = 1447.4166
$ws.Range("I135").Value = 1639.4
$ws.Range("J135").Value = 487.5
$ws.Range("K135").Value = 14754.6
$ws.Range("L135").Value = 4387.5
$ws.Range("M135").Value = -12219.6
$ws.Range("N135").Value = -9457.5

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H2").Value = 382.42856
$ws.Range("I2").Value = 462.9
$ws.Range("J2").Value = 181.25
$ws.Range("K2").Value = 462.9
$ws.Range("L2").Value = 181.25
$ws.Range("M2").Value = -349.9
$ws.Range("N2").Value = -407.25
$ws.Range("H23").Value = 756
$ws.Range("I23").Value = 12
$ws.Range("J23").Value = 1500
$ws.Range("K23").Value = 12
$ws.Range("L23").Value = 1500
$ws.Range("M23").Value = 211
$ws.Range("N23").Value = -1946
$ws.Range("H94").Value = 50000
$ws.Range("J94").Value = 50000
$ws.Range("L94").Value = 50000
$ws.Range("N94").Value = -51352
$ws.Range("H102").Value = 1485.25
$ws.Range("I102").Value = 1485.25
$ws.Range("K102").Value = 1485.25
$ws.Range("M102").Value = 136.75
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 5000
$ws.Range("K113").Value = 5000
$ws.Range("M113").Value = -2830

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H46").Value = 2550.75
$ws.Range("I46").Value = 2299.6155
$ws.Range("J46").Value = 3639
$ws.Range("K46").Value = 2299.6155
$ws.Range("L46").Value = 3639
$ws.Range("M46").Value = -2111.6155
$ws.Range("N46").Value = -4015
$ws.Range("H68").Value = 2169.2
$ws.Range("I68").Value = 2224
$ws.Range("J68").Value = 1950
$ws.Range("K68").Value = 2224
$ws.Range("L68").Value = 1950
$ws.Range("M68").Value = -1475
$ws.Range("N68").Value = -3448
$ws.Range("H71").Value = 2169.2
$ws.Range("I71").Value = 2224
$ws.Range("J71").Value = 1950
$ws.Range("K71").Value = 11120
$ws.Range("L71").Value = 9750
$ws.Range("M71").Value = -7376
$ws.Range("N71").Value = -17238

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 4746.7144
$ws.Range("J62").Value = 4365
$ws.Range("L62").Value = 4365
$ws.Range("N62").Value = -5613
$ws.Range("H65").Value = 4746.7144
$ws.Range("J65").Value = 4365
$ws.Range("N65").Value = -28065
$ws.Range("H124").Value = 23666.666
$ws.Range("J124").Value = 23666.666
$ws.Range("L124").Value = 23666.666
$ws.Range("N124").Value = -33486.666
